# DU FBS Mock 5 - results update
# Corrects raw MCQ correct/wrong counts for a handful of students; every
# dependent formula cell (marks, percentages, MCQ totals/ranks, grand
# totals/ranks) recalculates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Aymaan Zaman (row 6): Business Studies Correct/Wrong
$ws.Range("K6").Value = 7
$ws.Range("L6").Value = 3

# Mansib Rahman (row 9): Adv English Correct/Wrong
$ws.Range("G9").Value = 8
$ws.Range("H9").Value = 3

# Tarannum Rashid (row 12): Adv English + Business Studies Correct/Wrong
$ws.Range("G12").Value = 11
$ws.Range("H12").Value = 1
$ws.Range("K12").Value = 7
$ws.Range("L12").Value = 7

# Tasnia Taha (row 24): English Correct/Wrong
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 0

# Reflect the author's last on-screen selection before saving
[void]$ws.Range("AF12").Select()
